$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (Corequisites, Concurrent, Recommended) immediately
# before the existing "Terms Typically Offered" column, shifting that column
# from D to G and expanding the used range from A1:D29 to A1:G29.
$ws.Range("D1:F1").EntireColumn.Insert()

# New column headers for row 1
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default every new Corequisites / Concurrent / Recommended cell (rows 2-29) to "NA"
$ws.Range("D2:F29").Value = "NA"

# Row-specific overrides: course rows whose Prerequisites text contained an
# embedded "Recommended: ..." clause had that clause moved into the new
# Recommended column, and the Prerequisites text trimmed accordingly.

# Row 4: WGS 301. Contemporary Issues in Women's and Gender Studies.
$ws.Range("C4").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and two lower-division courses in GE Area D."
$ws.Range("F4").Value = "WGS 201 (GE Area D1)."

# Row 16: WGS 375. Intersectional Feminist Art Histories.
$ws.Range("C16").Value = "Junior standing; and one of the ART 112, ART 212, ART 213, or WGS 201."

# Row 18: WGS 401. Seminar in Women's and Gender Studies.
$ws.Range("C18").Value = "WGS 201."
$ws.Range("F18").Value = "Junior standing."

# Row 19: WGS 417. Feminist Legal Theory.
$ws.Range("C19").Value = "Completion of GE D1."
$ws.Range("F19").Value = "POLS 112."

# Row 26: WGS 457. U.S. Reproductive Politics.
$ws.Range("C26").Value = "Completion of GE D1."
$ws.Range("F26").Value = "POLS 112."

# A handful of "Terms Typically Offered" values carried a trailing space in the
# source data once they moved from column D to column G; reproduce that exactly.
$ws.Range("G4").Value = "F, W, SP "
$ws.Range("G18").Value = "W "
$ws.Range("G19").Value = "SP "
$ws.Range("G26").Value = "SP "

# Normalize a stray non-breaking space around course codes in these Prerequisites
# cells (e.g. "WGS<nbsp>201" -> "WGS 201") to match the rest of the sheet; the
# wording itself is unchanged for these rows.
$ws.Range("C9").Value = "PSY 201 or PSY 202."  # WGS 324. Psychology of Gender.
$ws.Range("C17").Value = "WGS 201 or consent of Women's and Gender Studies Chair."  # WGS 400. Special Problems for Advanced Undergraduates.
$ws.Range("C20").Value = "HIST 303; completion of GE Area D5; or graduate standing."  # WGS 421. The History of Prostitution.
$ws.Range("C22").Value = "HIST 303 or completion of GE Area D5, or graduate standing."  # WGS 434. American Women's History to 1870.
$ws.Range("C23").Value = "HIST 303 or completion of GE Area D5, or graduate standing."  # WGS 435. American Women's History from 1870.
$ws.Range("C24").Value = "WGS 201."  # WGS 450. Feminist Theory.
$ws.Range("C27").Value = "HIST 303 or completion of GE Area D5, or graduate standing."  # WGS 458. Gender and Sexuality in Modern Europe.
$ws.Range("C28").Value = "WGS 201 and consent of WGS/QS Internship Director."  # WGS 467. Women's and Gender Studies / Queer Studies Internship.
